$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Create Contact" (rows 6-11) and "Create Organization" (rows 12-16) blocks
# each end with a verification/result row whose Status column should read "Pass".
# Mark the Status cells for the "Create Contact" (row 7) and
# "Create Organization" (row 13) header/first-data rows as Pass.
$ws.Range("E7").Value = "Pass"
$ws.Range("E13").Value = "Pass"

# Insert a new row at row 16 (inside the "Create Organization" table), which
# pushes the existing terminating row ("####" / Pass row) down to row 17.
$ws.Rows.Item(16).Insert() | Out-Null

# Grow the third table (Table36, originally A12:E16) so it now covers the
# newly inserted row as well.
$lo = $ws.ListObjects.Item(3)
$lo.Resize($ws.Range("A12:E17")) | Out-Null

# Fill in the data for the freshly inserted row: a new "Last Name" test step.
$ws.Range("C16").Value = "New Last Name"
$ws.Range("D16").Value = "Cde"

# Leave the selection on the newly added cell, matching the edited workbook.
$ws.Range("D16").Select() | Out-Null
